# Adding Mark down for Regina's persona
# Remove the inline picture (Regina's portrait) from the document, leaving
# its paragraph empty.

$d = $word.ActiveDocument

if ($d.InlineShapes.Count -ge 1) {
    $d.InlineShapes.Item(1).Delete()
}
